# Update gh-pages output — refresh "want to go" counts, fix a start-time typo,
# and append the new "Nice Mini World" doll-expo row to the 展览 (Exhibition)
# and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value  = 620
$ws1.Range("F3").Value  = 583
$ws1.Range("F4").Value  = 884
$ws1.Range("F5").Value  = 632
$ws1.Range("F7").Value  = 363
$ws1.Range("F8").Value  = 554
$ws1.Range("F10").Value = 1127
$ws1.Range("F11").Value = 577
$ws1.Range("F12").Value = 338
$ws1.Range("F14").Value = 141
$ws1.Range("F15").Value = 294
$ws1.Range("F18").Value = 527
$ws1.Range("F19").Value = 16
$ws1.Range("F20").Value = 522
$ws1.Range("E21").Value = "2024.05.10 10:30-05.12 18:30"
$ws1.Range("F22").Value = 480

# New row 23 — copy column A's number style first so the index cell matches
# the sheet's existing "centered/bold/bordered" look, then fill the values.
$ws1.Range("A22").Copy()
$ws1.Range("A23").PasteSpecial(-4122)
$ws1.Range("A23").Value = 22

$ws1.Range("B23").NumberFormat = "@"
$ws1.Range("B23").Value = "2024.05.25"
$ws1.Range("B23").ClearFormats()

$ws1.Range("C23").Value = "广州·奶司的小人国娃展Nice Mini World  "
$ws1.Range("D23").Value = "洛浦街厦滘西环路1号 岭南会展中心"
$ws1.Range("E23").Value = "2024.05.25 10:30-05.25 17:00"
$ws1.Range("F23").Value = 0
$ws1.Range("G23").Value = 60
$ws1.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=82093"
$ws1.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202402/rhIj7fnH1708936497981.jpeg"

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performance)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F2").Value  = 68
$ws2.Range("F3").Value  = 54
$ws2.Range("F7").Value  = 636
$ws2.Range("F9").Value  = 195
$ws2.Range("F12").Value = 20

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) — same edits as sheet 1, mirrored at its own
# row offsets, plus its own copy of the new row (appended at the end, 35).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F3").Value  = 68
$ws4.Range("F4").Value  = 620
$ws4.Range("F5").Value  = 54
$ws4.Range("F7").Value  = 583
$ws4.Range("F8").Value  = 884
$ws4.Range("F9").Value  = 632
$ws4.Range("F11").Value = 363
$ws4.Range("F12").Value = 554
$ws4.Range("F14").Value = 1127
$ws4.Range("F15").Value = 577
$ws4.Range("F18").Value = 338
$ws4.Range("F20").Value = 636
$ws4.Range("F21").Value = 141
$ws4.Range("F23").Value = 294
$ws4.Range("F26").Value = 195
$ws4.Range("F28").Value = 527
$ws4.Range("F30").Value = 20
$ws4.Range("F31").Value = 16
$ws4.Range("F32").Value = 522
$ws4.Range("E33").Value = "2024.05.10 10:30-05.12 18:30"
$ws4.Range("F34").Value = 480

# New row 35 — identical content/index to sheet 1's new row 23.
$ws4.Range("A34").Copy()
$ws4.Range("A35").PasteSpecial(-4122)
$ws4.Range("A35").Value = 34

$ws4.Range("B35").NumberFormat = "@"
$ws4.Range("B35").Value = "2024.05.25"
$ws4.Range("B35").ClearFormats()

$ws4.Range("C35").Value = "广州·奶司的小人国娃展Nice Mini World  "
$ws4.Range("D35").Value = "洛浦街厦滘西环路1号 岭南会展中心"
$ws4.Range("E35").Value = "2024.05.25 10:30-05.25 17:00"
$ws4.Range("F35").Value = 0
$ws4.Range("G35").Value = 60
$ws4.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=82093"
$ws4.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202402/rhIj7fnH1708936497981.jpeg"
